$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D11 (row 11, 07.02) already carries the wrapped-text/bordered "filled-in
# row" style (s=10); reuse that formatting for the newly-filled E11 (comment)
# and D12 (goal, row 12 = 07.05) cells before writing their text, so the new
# cells pick up the same look as the rest of the table instead of the
# blank-row default style.
$ws.Range("D11").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Row 11 (07.02): update the Goal text (3rd bullet changed) and fill in the
# previously-empty Comment cell.
$ws.Range("D11").Value = "1. Git 목록 삭제`n2. 요구사항 분석`n3. 인터페이스 구상도 작성"
$ws.Range("E11").Value = "1. 07.02 회의록 작성`n2. Git 업로드 한 것 중 중복된 것들 삭제`n3. 어플 흐름 구상도 작성`n4. 요구사항/기능 정의/기능 설계 3가지 부분으로 나누어 요구사항 수정(미완)"

# Row 12 (07.05): fill in the previously-empty Goal cell.
$ws.Range("D12").Value = "1. 인터페이스 구상도 작성`n2. 요구사항 분류하여 수정"

# The extra lines of text mean these two rows now need to be taller.
$ws.Rows(11).RowHeight = 65.85
$ws.Rows(12).RowHeight = 26.35

# Move the active selection, matching where the author ended up editing.
$ws.Range("D4").Select() | Out-Null
